# Inserts one new weekly price record for "Betarraga" (Terminal
# Hortofrutícola Agro Chillán) immediately above the current row 193,
# pushing the existing rows 193:291 down to 194:292 (dimension A1:R291 ->
# A1:R292).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 193 and below down by one row.
$ws.Rows.Item(193).Insert()

# Populate the newly opened row 193 with the new record.
$ws.Cells.Item(193, 1).Value  = 7
$ws.Cells.Item(193, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(193, 3).Value  = "Ñuble"
$ws.Cells.Item(193, 4).Value  = 44529
$ws.Cells.Item(193, 5).Value  = 16
$ws.Cells.Item(193, 6).Value  = 100114014
$ws.Cells.Item(193, 7).Value  = "Betarraga"
$ws.Cells.Item(193, 8).Value  = "Sin especificar"
$ws.Cells.Item(193, 9).Value  = "Primera"
$ws.Cells.Item(193, 10).Value = 240
$ws.Cells.Item(193, 11).Value = 700
$ws.Cells.Item(193, 12).Value = 800
$ws.Cells.Item(193, 13).Value = 750
$ws.Cells.Item(193, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(193, 15).Value = "Región del Maule"
$ws.Cells.Item(193, 16).Value = 150
$ws.Cells.Item(193, 17).Value = 5
$ws.Cells.Item(193, 18).Value = "Hortaliza"
